# PASO 18 modifico en rama htmlify
# Wrap each verse of the poem in <p>/<em>/<code> "HTML-ified" markup,
# rewrite run text to the escaped HTML strings, change "git" -> "pull"
# in the fourth verse, and tidy up the trailing empty paragraphs.

$d = $word.ActiveDocument

$wdFindContinue = 1
$wdReplaceAll = 2

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
        $true, $wdFindContinue, $false, $replace, $wdReplaceAll) | Out-Null
}

# --- Verses 1-4: single/multi run paragraphs where Find/Replace already
#     collapses every run (and any *interior* proofErr pair) into one run. ---

Replace-Text "En un lugar de los Repos," "<p>En un lugar de los <em>Repos</em>,</p>"

Replace-Text "de cuyo commit no quiero acordarme," "<p>de cuyo <em>commit</em> no quiero acordarme,</p>"

Replace-Text "no ha mucho log que vivía" "<p>no ha mucho <em>log</em> que vivía</p>"

Replace-Text "un hidalgo de los de head en master," "<p>un hidalgo de los de <em>head</em> en <em>master</em>,</p>"

# --- Verses 5-7: each paragraph starts with a lone <w:proofErr spellStart/>
#     immediately after <w:pPr> (i.e. before the first run). Find/Replace
#     happily merges/collapses interior runs (and the proofErr pairs that
#     sit between them), but it leaves that *leading* marker behind since
#     it isn't "inside" any run. Rebuild those three paragraphs from
#     scratch (insert a fresh paragraph, fill it in, drop the old one) so
#     the stray leading proofErr goes away too. ---

function Rebuild-Paragraph($index, $newText) {
    $old = $d.Paragraphs($index)
    $old.Range.InsertParagraphAfter() | Out-Null
    $new = $d.Paragraphs($index + 1)
    $new.Range.Text = $newText
    $d.Paragraphs($index).Range.Delete() | Out-Null
}

Rebuild-Paragraph 5 "<p><em>push</em> antiguo,</p>"

Rebuild-Paragraph 6 "<p><em>pull</em> flaco y <em>remote</em> corredor.</p>"

Rebuild-Paragraph 7 "<p><code>git commit -m &quot;Don Quijote de la Mancha&quot;</code> </p>"

# --- Trailing paragraphs: the empty "Times New Roman, 20half-pt" paragraph
#     right before the _GoBack bookmark paragraph goes away, its paragraph
#     mark run properties move onto the bookmark paragraph, and a brand
#     new empty paragraph is appended after it (before the sectPr). ---

$blankIndex = $d.Paragraphs.Count - 1
$bookmarkIndex = $d.Paragraphs.Count

$blank = $d.Paragraphs($blankIndex)
$bookmark = $d.Paragraphs($bookmarkIndex)

$bookmark.Range.Font.Name = "Times"
$bookmark.Range.Font.Size = 10

$blank.Range.Delete() | Out-Null

$d.Paragraphs($d.Paragraphs.Count).Range.InsertParagraphAfter() | Out-Null
